$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet to reflect new "through" date
$ws.Name = "Through 2022-02-11"

# Update the "2022 (through ...)" header label (column I header, row 1)
$ws.Range("I1").Value = "2022 (through 02-11)"

# Update affected monthly/annual carjacking totals for the 2022 column (I)
$ws.Range("I2").Value = 161   # January
$ws.Range("I3").Value = 52    # February
$ws.Range("I14").Value = 213  # Total
